$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant columns shared across all data rows in this dataset
$constA = 10
$constB = "Vega Modelo de Temuco"
$constC = "La Araucanía"
$constE = 9
$constF = "Fruta"
$constG = 100103
$constH = "Frutos de hueso (carozo)"
$constI = 100103001
$constJ = "Cereza"

# Ensure constant columns + date format are populated for newly added rows 125-127
foreach ($r in 125..127) {
    $ws.Cells.Item($r, 1).Value = $constA
    $ws.Cells.Item($r, 2).Value = $constB
    $ws.Cells.Item($r, 3).Value = $constC
    $ws.Cells.Item($r, 5).Value = $constE
    $ws.Cells.Item($r, 6).Value = $constF
    $ws.Cells.Item($r, 7).Value = $constG
    $ws.Cells.Item($r, 8).Value = $constH
    $ws.Cells.Item($r, 9).Value = $constI
    $ws.Cells.Item($r, 10).Value = $constJ
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# Update D, K, L, M, N, O, P, Q, R, S, T for rows 105-127 to reflect the new weekly data
$ws.Cells.Item(105, 4).Value = 44522
$ws.Cells.Item(105, 11).Value = "Early Burlat"
$ws.Cells.Item(105, 12).Value = "Especial"
$ws.Cells.Item(105, 13).Value = 100
$ws.Cells.Item(105, 14).Value = 2600
$ws.Cells.Item(105, 15).Value = 2600
$ws.Cells.Item(105, 16).Value = 2600
$ws.Cells.Item(105, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(105, 18).Value = "Región del Maule"
$ws.Cells.Item(105, 19).Value = 2600
$ws.Cells.Item(105, 20).Value = 1

$ws.Cells.Item(106, 4).Value = 44522
$ws.Cells.Item(106, 11).Value = "Early Burlat"
$ws.Cells.Item(106, 12).Value = "Primera"
$ws.Cells.Item(106, 13).Value = 50
$ws.Cells.Item(106, 14).Value = 20000
$ws.Cells.Item(106, 15).Value = 20000
$ws.Cells.Item(106, 16).Value = 20000
$ws.Cells.Item(106, 17).Value = "`$/caja 10 kilos"
$ws.Cells.Item(106, 18).Value = "Región del Maule"
$ws.Cells.Item(106, 19).Value = 2000
$ws.Cells.Item(106, 20).Value = 10

$ws.Cells.Item(107, 4).Value = 44522
$ws.Cells.Item(107, 11).Value = "Early Burlat"
$ws.Cells.Item(107, 12).Value = "Primera"
$ws.Cells.Item(107, 13).Value = 400
$ws.Cells.Item(107, 14).Value = 2000
$ws.Cells.Item(107, 15).Value = 2200
$ws.Cells.Item(107, 16).Value = 2100
$ws.Cells.Item(107, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(107, 18).Value = "Región del Maule"
$ws.Cells.Item(107, 19).Value = 2100
$ws.Cells.Item(107, 20).Value = 1

$ws.Cells.Item(108, 4).Value = 44166
$ws.Cells.Item(108, 11).Value = "Brooks"
$ws.Cells.Item(108, 12).Value = "Primera"
$ws.Cells.Item(108, 13).Value = 200
$ws.Cells.Item(108, 14).Value = 2000
$ws.Cells.Item(108, 15).Value = 2000
$ws.Cells.Item(108, 16).Value = 2000
$ws.Cells.Item(108, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(108, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(108, 19).Value = 2000
$ws.Cells.Item(108, 20).Value = 1

$ws.Cells.Item(109, 4).Value = 44166
$ws.Cells.Item(109, 11).Value = "Corazón de paloma"
$ws.Cells.Item(109, 12).Value = "Primera"
$ws.Cells.Item(109, 13).Value = 140
$ws.Cells.Item(109, 14).Value = 2000
$ws.Cells.Item(109, 15).Value = 2000
$ws.Cells.Item(109, 16).Value = 2000
$ws.Cells.Item(109, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(109, 18).Value = "Región del Maule"
$ws.Cells.Item(109, 19).Value = 2000
$ws.Cells.Item(109, 20).Value = 1

$ws.Cells.Item(110, 4).Value = 44166
$ws.Cells.Item(110, 11).Value = "Santina"
$ws.Cells.Item(110, 12).Value = "Primera"
$ws.Cells.Item(110, 13).Value = 300
$ws.Cells.Item(110, 14).Value = 1500
$ws.Cells.Item(110, 15).Value = 1500
$ws.Cells.Item(110, 16).Value = 1500
$ws.Cells.Item(110, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(110, 18).Value = "Región del Maule"
$ws.Cells.Item(110, 19).Value = 1500
$ws.Cells.Item(110, 20).Value = 1

$ws.Cells.Item(111, 4).Value = 44232
$ws.Cells.Item(111, 11).Value = "Bing"
$ws.Cells.Item(111, 12).Value = "Primera"
$ws.Cells.Item(111, 13).Value = 185
$ws.Cells.Item(111, 14).Value = 900
$ws.Cells.Item(111, 15).Value = 1000
$ws.Cells.Item(111, 16).Value = 951
$ws.Cells.Item(111, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(111, 18).Value = "Región del Maule"
$ws.Cells.Item(111, 19).Value = 951
$ws.Cells.Item(111, 20).Value = 1

$ws.Cells.Item(112, 4).Value = 44161
$ws.Cells.Item(112, 11).Value = "Royal Dawn"
$ws.Cells.Item(112, 12).Value = "Especial"
$ws.Cells.Item(112, 13).Value = 155
$ws.Cells.Item(112, 14).Value = 2100
$ws.Cells.Item(112, 15).Value = 2100
$ws.Cells.Item(112, 16).Value = 2100
$ws.Cells.Item(112, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(112, 18).Value = "Región del Maule"
$ws.Cells.Item(112, 19).Value = 2100
$ws.Cells.Item(112, 20).Value = 1

$ws.Cells.Item(113, 4).Value = 44161
$ws.Cells.Item(113, 11).Value = "Royal Dawn"
$ws.Cells.Item(113, 12).Value = "Primera"
$ws.Cells.Item(113, 13).Value = 210
$ws.Cells.Item(113, 14).Value = 2000
$ws.Cells.Item(113, 15).Value = 2000
$ws.Cells.Item(113, 16).Value = 2000
$ws.Cells.Item(113, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(113, 18).Value = "Región del Maule"
$ws.Cells.Item(113, 19).Value = 2000
$ws.Cells.Item(113, 20).Value = 1

$ws.Cells.Item(114, 4).Value = 44161
$ws.Cells.Item(114, 11).Value = "Royal Dawn"
$ws.Cells.Item(114, 12).Value = "Segunda"
$ws.Cells.Item(114, 13).Value = 260
$ws.Cells.Item(114, 14).Value = 1500
$ws.Cells.Item(114, 15).Value = 1500
$ws.Cells.Item(114, 16).Value = 1500
$ws.Cells.Item(114, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(114, 18).Value = "Región del Maule"
$ws.Cells.Item(114, 19).Value = 1500
$ws.Cells.Item(114, 20).Value = 1

$ws.Cells.Item(115, 4).Value = 44238
$ws.Cells.Item(115, 11).Value = "Lapins"
$ws.Cells.Item(115, 12).Value = "Primera"
$ws.Cells.Item(115, 13).Value = 225
$ws.Cells.Item(115, 14).Value = 1000
$ws.Cells.Item(115, 15).Value = 1200
$ws.Cells.Item(115, 16).Value = 1111
$ws.Cells.Item(115, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(115, 18).Value = "Región del Maule"
$ws.Cells.Item(115, 19).Value = 1111
$ws.Cells.Item(115, 20).Value = 1

$ws.Cells.Item(116, 4).Value = 44515
$ws.Cells.Item(116, 11).Value = "Early Burlat"
$ws.Cells.Item(116, 12).Value = "Primera"
$ws.Cells.Item(116, 13).Value = 30
$ws.Cells.Item(116, 14).Value = 3500
$ws.Cells.Item(116, 15).Value = 3500
$ws.Cells.Item(116, 16).Value = 3500
$ws.Cells.Item(116, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(116, 18).Value = "Región del Maule"
$ws.Cells.Item(116, 19).Value = 3500
$ws.Cells.Item(116, 20).Value = 1

$ws.Cells.Item(117, 4).Value = 44181
$ws.Cells.Item(117, 11).Value = "Lapins"
$ws.Cells.Item(117, 12).Value = "Primera"
$ws.Cells.Item(117, 13).Value = 300
$ws.Cells.Item(117, 14).Value = 1000
$ws.Cells.Item(117, 15).Value = 1000
$ws.Cells.Item(117, 16).Value = 1000
$ws.Cells.Item(117, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(117, 18).Value = "Región del Maule"
$ws.Cells.Item(117, 19).Value = 1000
$ws.Cells.Item(117, 20).Value = 1

$ws.Cells.Item(118, 4).Value = 44181
$ws.Cells.Item(118, 11).Value = "Rainier"
$ws.Cells.Item(118, 12).Value = "Primera"
$ws.Cells.Item(118, 13).Value = 200
$ws.Cells.Item(118, 14).Value = 1000
$ws.Cells.Item(118, 15).Value = 1200
$ws.Cells.Item(118, 16).Value = 1100
$ws.Cells.Item(118, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(118, 18).Value = "Región del Maule"
$ws.Cells.Item(118, 19).Value = 1100
$ws.Cells.Item(118, 20).Value = 1

$ws.Cells.Item(119, 4).Value = 44181
$ws.Cells.Item(119, 11).Value = "Royal Dawn"
$ws.Cells.Item(119, 12).Value = "Primera"
$ws.Cells.Item(119, 13).Value = 700
$ws.Cells.Item(119, 14).Value = 1000
$ws.Cells.Item(119, 15).Value = 1000
$ws.Cells.Item(119, 16).Value = 1000
$ws.Cells.Item(119, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(119, 18).Value = "Región del Maule"
$ws.Cells.Item(119, 19).Value = 1000
$ws.Cells.Item(119, 20).Value = 1

$ws.Cells.Item(120, 4).Value = 44181
$ws.Cells.Item(120, 11).Value = "Royal Dawn"
$ws.Cells.Item(120, 12).Value = "Segunda"
$ws.Cells.Item(120, 13).Value = 100
$ws.Cells.Item(120, 14).Value = 800
$ws.Cells.Item(120, 15).Value = 800
$ws.Cells.Item(120, 16).Value = 800
$ws.Cells.Item(120, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(120, 18).Value = "Región del Maule"
$ws.Cells.Item(120, 19).Value = 800
$ws.Cells.Item(120, 20).Value = 1

$ws.Cells.Item(121, 4).Value = 44181
$ws.Cells.Item(121, 11).Value = "Santina"
$ws.Cells.Item(121, 12).Value = "Primera"
$ws.Cells.Item(121, 13).Value = 200
$ws.Cells.Item(121, 14).Value = 1000
$ws.Cells.Item(121, 15).Value = 1000
$ws.Cells.Item(121, 16).Value = 1000
$ws.Cells.Item(121, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(121, 18).Value = "Región del Maule"
$ws.Cells.Item(121, 19).Value = 1000
$ws.Cells.Item(121, 20).Value = 1

$ws.Cells.Item(122, 4).Value = 44194
$ws.Cells.Item(122, 11).Value = "Brooks"
$ws.Cells.Item(122, 12).Value = "Primera"
$ws.Cells.Item(122, 13).Value = 500
$ws.Cells.Item(122, 14).Value = 10000
$ws.Cells.Item(122, 15).Value = 11000
$ws.Cells.Item(122, 16).Value = 10600
$ws.Cells.Item(122, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(122, 18).Value = "Región del Maule"
$ws.Cells.Item(122, 19).Value = 1060
$ws.Cells.Item(122, 20).Value = 10

$ws.Cells.Item(123, 4).Value = 44194
$ws.Cells.Item(123, 11).Value = "Brooks"
$ws.Cells.Item(123, 12).Value = "Primera"
$ws.Cells.Item(123, 13).Value = 220
$ws.Cells.Item(123, 14).Value = 600
$ws.Cells.Item(123, 15).Value = 700
$ws.Cells.Item(123, 16).Value = 645
$ws.Cells.Item(123, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(123, 18).Value = "Región del Maule"
$ws.Cells.Item(123, 19).Value = 645
$ws.Cells.Item(123, 20).Value = 1

$ws.Cells.Item(124, 4).Value = 44194
$ws.Cells.Item(124, 11).Value = "Corazón de paloma"
$ws.Cells.Item(124, 12).Value = "Especial"
$ws.Cells.Item(124, 13).Value = 100
$ws.Cells.Item(124, 14).Value = 14000
$ws.Cells.Item(124, 15).Value = 14000
$ws.Cells.Item(124, 16).Value = 14000
$ws.Cells.Item(124, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(124, 18).Value = "Región del Maule"
$ws.Cells.Item(124, 19).Value = 1400
$ws.Cells.Item(124, 20).Value = 10

$ws.Cells.Item(125, 4).Value = 44194
$ws.Cells.Item(125, 11).Value = "Lapins"
$ws.Cells.Item(125, 12).Value = "Primera"
$ws.Cells.Item(125, 13).Value = 1100
$ws.Cells.Item(125, 14).Value = 600
$ws.Cells.Item(125, 15).Value = 700
$ws.Cells.Item(125, 16).Value = 655
$ws.Cells.Item(125, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(125, 18).Value = "Región del Maule"
$ws.Cells.Item(125, 19).Value = 655
$ws.Cells.Item(125, 20).Value = 1

$ws.Cells.Item(126, 4).Value = 44201
$ws.Cells.Item(126, 11).Value = "Brooks"
$ws.Cells.Item(126, 12).Value = "Primera"
$ws.Cells.Item(126, 13).Value = 235
$ws.Cells.Item(126, 14).Value = 900
$ws.Cells.Item(126, 15).Value = 1000
$ws.Cells.Item(126, 16).Value = 953
$ws.Cells.Item(126, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(126, 18).Value = "Región del Maule"
$ws.Cells.Item(126, 19).Value = 953
$ws.Cells.Item(126, 20).Value = 1

$ws.Cells.Item(127, 4).Value = 44201
$ws.Cells.Item(127, 11).Value = "Lapins"
$ws.Cells.Item(127, 12).Value = "Primera"
$ws.Cells.Item(127, 13).Value = 390
$ws.Cells.Item(127, 14).Value = 800
$ws.Cells.Item(127, 15).Value = 900
$ws.Cells.Item(127, 16).Value = 846
$ws.Cells.Item(127, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(127, 18).Value = "Región del Maule"
$ws.Cells.Item(127, 19).Value = 846
$ws.Cells.Item(127, 20).Value = 1
